$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated crypto price/volume data (GitHub Actions refresh).
# Column D (Price) values are prefixed with a leading apostrophe so Excel
# stores them as literal text (matching the original inline-string cells)
# instead of auto-converting number-looking strings (e.g. "1.00", "0.424")
# into numeric values.

$ws.Range("D2").Value = "'96.996.83"
$ws.Range("E2").Value = '  +0.62%  '
$ws.Range("D3").Value = "'3.694.92"
$ws.Range("E3").Value = '  +0.94%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'237.14"
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("E6").Value = '  +0.76%  '
$ws.Range("D7").Value = "'658.69"
$ws.Range("E7").Value = '  -0.23%  '
$ws.Range("D8").Value = "'0.424"
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -2.08%  '
$ws.Range("D11").Value = "'3.694.14"
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").Value = "'44.29"
$ws.Range("E12").Value = '  -1.42%  '
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("D14").Value = "'0.0000300"
$ws.Range("E14").Value = '  +11.18%  '
$ws.Range("E15").Value = '  +1.57%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = "'4.374.96"
$ws.Range("E16").Value = '  +0.80%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = "'96.835.15"
$ws.Range("E17").Value = '  +0.73%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = "'9.18"
$ws.Range("E18").Value = '  +3.05%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = "'3.703.28"
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = "'13.04"
$ws.Range("E20").Value = '  +2.47%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = "'18.78"
$ws.Range("E21").Value = '  +2.96%  '
$ws.Range("B22").Value = 'Stellar'
$ws.Range("C22").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D22").Value = "'0.508"
$ws.Range("E22").Value = '  -4.87%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = "'519.67"
$ws.Range("E23").Value = '  -0.23%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = "'3.42"
$ws.Range("E24").Value = '  -0.80%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = "'0.0000210"
$ws.Range("E25").Value = '  +3.16%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = "'6.95"
$ws.Range("E26").Value = '  +0.79%  '
$ws.Range("B27").Value = 'Hedera'
$ws.Range("C27").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D27").Value = "'0.204"
$ws.Range("E27").Value = '  +23.46%  '
$ws.Range("B28").Value = 'Litecoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D28").Value = "'101.32"
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = "'13.39"
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'12.55"
$ws.Range("E30").Value = '  +2.13%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'3.02"
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("B32").Value = 'Dai'
$ws.Range("C32").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").Value = "'0.192"
$ws.Range("E33").Value = '  +3.90%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").Value = "'1.87"
$ws.Range("E34").Value = '  +2.36%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = '  -0.50%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = "'32.19"
$ws.Range("E36").Value = '  -2.46%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = "'646.59"
$ws.Range("E37").Value = '  +3.38%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").Value = "'0.593"
$ws.Range("E38").Value = '  +1.29%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = "'8.86"
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = "'1.00"
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = "'0.504"
$ws.Range("E41").Value = '  +16.25%  '
$ws.Range("D42").Value = "'2.08"
$ws.Range("E42").Value = '  +6.85%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = "'6.88"
$ws.Range("E43").Value = '  +9.43%  '
$ws.Range("B44").Value = 'Kaspa'
$ws.Range("C44").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D44").Value = "'0.161"
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = "'40.53"
$ws.Range("E45").Value = '  -10.37%  '
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = "'0.962"
$ws.Range("E46").Value = '  +0.57%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = "'0.0468"
$ws.Range("E47").Value = '  +3.10%  '
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D49").Value = "'23.63"
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = "'8.67"
$ws.Range("E50").Value = '  +1.67%  '
$ws.Range("B51").Value = 'MantraDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D51").Value = "'3.53"
$ws.Range("E51").Value = '  -1.99%  '
